# Weekly update: insert a new daily record for "Hass / Primera" at row 58
# (Perú origin, $/bandeja 10 kilos), pushing all subsequent rows (old 58-99)
# down by one row (new 59-100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 58; existing rows 58-99 shift to 59-100.
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new record.
$ws.Range("A58").Value = 1
$ws.Range("B58").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C58").Value = "Arica y Parinacota"
$ws.Range("D58").Value = 44729
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100106
$ws.Range("H58").Value = "Oleaginosos"
$ws.Range("I58").Value = 100106002
$ws.Range("J58").Value = "Palta"
$ws.Range("K58").Value = "Hass"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 400
$ws.Range("N58").Value = 16000
$ws.Range("O58").Value = 17000
$ws.Range("P58").Value = 16500
$ws.Range("Q58").Value = '$/bandeja 10 kilos'
$ws.Range("R58").Value = "Perú"
$ws.Range("S58").Value = 1650
$ws.Range("T58").Value = 10
